# Update cryptocurrency price (D) and 1h volume-change (E) columns
# Values are forced to Text via a leading quote-prefix so that
# number-looking strings (e.g. "0.9990", "2.050") keep their exact
# textual representation instead of being parsed as numbers, then the
# style is reset to "Normal" so the quote-prefix flag does not leave a
# residual cell style (matching the original inline-string cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.307.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "'1.811.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.86%  "
$ws.Range("D5").Value = "'338.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'0.9958"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").Value = "'0.4372"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +15.13%  "
$ws.Range("D8").Value = "'0.3531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.08%  "
$ws.Range("D9").Value = "'45.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "'1.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").Value = "'0.07459"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "'22.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "'0.9975"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'6.291"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'7.294"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "'1.814.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "'0.06672"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'82.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "'0.9969"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'6.458"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").Value = "'28.318.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").Value = "'2.385"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'2.468"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.36%  "
$ws.Range("D27").Value = "'20.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("D28").Value = "'155.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").Value = "'2.024.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").Value = "'1.311"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.24%  "
$ws.Range("D31").Value = "'132.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").Value = "'4.064"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "'5.978"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").Value = "'0.09324"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("D35").Value = "'12.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "'0.6821"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'0.02383"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "'0.06273"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").Value = "'5.219"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'0.2166"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").Value = "'1.494"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'1.221"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "'8.265"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").Value = "'0.9953"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "'13.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "'0.6181"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'3.869"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "'129.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("D49").Value = "'2.050"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'1.176"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "'0.07114"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.08%  "
